$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The caption paragraph currently reads:   "ตาราง … Activity Diagram"
# and must become:                          "ตารางที่ 1 Activity Diagram"
#
# i.e.
#   run "ตาราง"            -> unchanged
#   run " "                -> "ที่ "
#   run "… "                -> split into "1" and " "
#   run "Activity Diagram" -> unchanged
# ------------------------------------------------------------------

# Locate the three anchor runs by text (robust to any offset changes).
$rngWord = $d.Content.Duplicate
[void]$rngWord.Find.Execute("ตาราง", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wordStart = $rngWord.Start
$wordEnd   = $rngWord.End

$rngEllipsis = $d.Content.Duplicate
[void]$rngEllipsis.Find.Execute("… ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ellipsisStart = $rngEllipsis.Start
$ellipsisEnd   = $rngEllipsis.End

$rngTail = $d.Content.Duplicate
[void]$rngTail.Find.Execute("Activity Diagram", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailStart = $rngTail.Start
$tailEnd   = $rngTail.End

# ------------------------------------------------------------------
# Phase 1: temporarily mark (Bold, used only as a scratch flag) the
# runs that must stay distinct from a soon-to-be-identically
# formatted neighbour, so the text edits below do not get silently
# coalesced into that neighbour.
# ------------------------------------------------------------------
$d.Range($wordStart, $wordEnd).Bold = 1
$d.Range($tailStart, $tailEnd).Bold = 1

# ------------------------------------------------------------------
# Phase 2: perform the textual substitutions, right-to-left so that
# the previously located offsets stay valid.
#   "… "  -> "1 "   (will be split into "1" + " " in phase 3)
#   " "   -> "ที่ "  (the space that directly follows "ตาราง")
# ------------------------------------------------------------------
$d.Range($ellipsisStart, $ellipsisEnd).Text = "1 "
$d.Range($wordEnd, $ellipsisStart).Text = "ที่ "

# ------------------------------------------------------------------
# Phase 3: split "1 " into two runs ("1" and " ") by marking the
# space with the same scratch marker, forcing it to stay separate
# from the "1" run that precedes it.
# ------------------------------------------------------------------
$rngNum = $d.Content.Duplicate
[void]$rngNum.Find.Execute("1 Activity Diagram", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$numStart = $rngNum.Start
$d.Range($numStart + 1, $numStart + 2).Bold = 1

# ------------------------------------------------------------------
# Phase 4: clear every scratch marker now that all runs that needed
# to stay separate already exist as distinct runs. Pure formatting
# changes do not trigger the run-coalescing pass, so this is safe.
# ------------------------------------------------------------------
$d.Range($wordStart, $wordEnd).Bold = 0
$d.Range($numStart + 1, $numStart + 2).Bold = 0

$rngTail2 = $d.Content.Duplicate
[void]$rngTail2.Find.Execute("Activity Diagram", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngTail2.Bold = 0
